$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Invoice")

# Row 8: invoice number 59 -> 60
$ws.Range("D8").Value = 60

# Row 18: remove the "Juustokakku 12 hlö" line item (product, qty, price) -
# clear the content but keep the row's formatting/formula in column E
$ws.Range("B18:D18").ClearContents()

# Update selection / scroll position to match the saved view
$ws.Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
$ws.Range("L11").Select()
